$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cdm-contact-point"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Remove the compound "cpt-2 / ele-1" constraint text from the root ContactPoint row's
# Constraint(s) column (AI2), leaving it empty - matching the removed shared string.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
